# Update the cryptos price list (Price / Volume(1h) columns) with the latest
# scraped values from coinranking.com, as produced by the scheduled
# GitHub Actions job. Row 47/48 also swap which coin (PancakeSwap /
# NEARProtocol) occupies which row, together with their refreshed values.
#
# Numeric-looking text values in column D (e.g. "1.017", "27.979.32") must
# stay stored as text, exactly as in the source data, so we prefix them
# with a leading apostrophe to force Excel to treat them as text rather
# than auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.979.32"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "'1.888.73"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "'1.017"
$ws.Range("E4").Value = "  +1.42%  "
$ws.Range("D5").Value = "'335.94"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("D7").Value = "'0.4697"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D8").Value = "'0.3939"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").Value = "'46.87"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "'0.08000"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "'21.80"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "'1.889.86"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "'5.984"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "'7.166"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "'0.06765"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "'87.81"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'0.00001051"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'17.20"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").Value = "'1.015"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").Value = "'27.983.96"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "'2.363"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").Value = "'2.114.16"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'158.98"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "'20.04"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'2.105"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'5.501"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").Value = "'121.56"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'0.09576"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'0.9658"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("D34").Value = "'3.650"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'5.361"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").Value = "'1.362"
$ws.Range("E36").Value = "  -7.03%  "
$ws.Range("D37").Value = "'0.06133"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").Value = "'1.214"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "'8.224"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'0.5969"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'0.1903"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'10.34"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'1.265"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "'0.5697"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "'12.23"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.947"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.397"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "'0.06873"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'113.75"
$ws.Range("D51").Value = "'1.070"
